# Auto-generated cell updates derived from the canonical OOXML diff.
# Applies new numeric values to Statistic sheet rows 3-6 (per-configuration stats).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("AB3").Value = 0.7222222222222222
$ws.Range("AC3").Value = 33200
$ws.Range("AD3").Value = 16
$ws.Range("AH3").Value = 0.3125
$ws.Range("AK3").Value = 0.6875
$ws.Range("AL3").Value = 36600
$ws.Range("AM3").Value = 14
$ws.Range("AQ3").Value = 0.35714285714285715
$ws.Range("AT3").Value = 0.6428571428571429
$ws.Range("AU3").Value = 43000
$ws.Range("AV3").Value = 8
$ws.Range("AZ3").Value = 0.625
$ws.Range("B3").Value = 31000
$ws.Range("BC3").Value = 0.375
$ws.Range("BD3").Value = 36500
$ws.Range("BE3").Value = 13
$ws.Range("BI3").Value = 0.38461538461538464
$ws.Range("BL3").Value = 0.6153846153846154
$ws.Range("BM3").Value = 36100
$ws.Range("BN3").Value = 14
$ws.Range("BR3").Value = 0.35714285714285715
$ws.Range("BU3").Value = 0.6428571428571429
$ws.Range("BV3").Value = 34600
$ws.Range("BW3").Value = 15
$ws.Range("C3").Value = 14
$ws.Range("CA3").Value = 0.3333333333333333
$ws.Range("CD3").Value = 0.6666666666666666
$ws.Range("CE3").Value = 29900
$ws.Range("CF3").Value = 16
$ws.Range("CG3").Value = 5
$ws.Range("CH3").Value = 5
$ws.Range("CJ3").Value = 0.3125
$ws.Range("CK3").Value = 1
$ws.Range("CM3").Value = 0.6875
$ws.Range("CN3").Value = 33970
$ws.Range("CO3").Value = 14.5
$ws.Range("CP3").Value = 5
$ws.Range("CQ3").Value = 5
$ws.Range("CS3").Value = 0.36112727142138906
$ws.Range("CT3").Value = 1
$ws.Range("CV3").Value = 0.6388727285786111
$ws.Range("CW3").Value = 28300
$ws.Range("CX3").Value = 18
$ws.Range("DB3").Value = 0.625
$ws.Range("DE3").Value = 0.375
$ws.Range("G3").Value = 0.35714285714285715
$ws.Range("J3").Value = 0.6428571428571429
$ws.Range("K3").Value = 30500
$ws.Range("L3").Value = 17
$ws.Range("M3").Value = 5
$ws.Range("N3").Value = 5
$ws.Range("P3").Value = 0.29411764705882354
$ws.Range("Q3").Value = 1
$ws.Range("S3").Value = 0.7058823529411765
$ws.Range("T3").Value = 28300
$ws.Range("U3").Value = 18
$ws.Range("Y3").Value = 0.2777777777777778

# Row 4
$ws.Range("AB4").Value = 0.868421052631579
$ws.Range("AC4").Value = 37700
$ws.Range("AD4").Value = 44
$ws.Range("AH4").Value = 0.11363636363636363
$ws.Range("AK4").Value = 0.8863636363636364
$ws.Range("AL4").Value = 30800
$ws.Range("AM4").Value = 38
$ws.Range("AQ4").Value = 0.13157894736842105
$ws.Range("AT4").Value = 0.868421052631579
$ws.Range("AU4").Value = 16400
$ws.Range("AV4").Value = 37
$ws.Range("AZ4").Value = 0.13513513513513514
$ws.Range("B4").Value = 46300
$ws.Range("BC4").Value = 0.8648648648648649
$ws.Range("BD4").Value = 12200
$ws.Range("BE4").Value = 49
$ws.Range("BI4").Value = 0.10204081632653061
$ws.Range("BL4").Value = 0.8979591836734694
$ws.Range("BM4").Value = 36000
$ws.Range("BN4").Value = 37
$ws.Range("BR4").Value = 0.13513513513513514
$ws.Range("BU4").Value = 0.8648648648648649
$ws.Range("BV4").Value = 13000
$ws.Range("BW4").Value = 48
$ws.Range("C4").Value = 42
$ws.Range("CA4").Value = 0.10416666666666667
$ws.Range("CD4").Value = 0.8958333333333334
$ws.Range("CE4").Value = 12700
$ws.Range("CF4").Value = 45
$ws.Range("CG4").Value = 5
$ws.Range("CH4").Value = 5
$ws.Range("CJ4").Value = 0.1111111111111111
$ws.Range("CK4").Value = 1
$ws.Range("CM4").Value = 0.8888888888888888
$ws.Range("CN4").Value = 27470
$ws.Range("CO4").Value = 41.1
$ws.Range("CP4").Value = 4.9
$ws.Range("CQ4").Value = 4.9
$ws.Range("CS4").Value = 0.12111363695010313
$ws.Range("CT4").Value = 0.9800000000000001
$ws.Range("CV4").Value = 0.878886363049897
$ws.Range("CW4").Value = 12200
$ws.Range("CX4").Value = 49
$ws.Range("D4").Value = 4
$ws.Range("DB4").Value = 0.15151515151515152
$ws.Range("DE4").Value = 0.8484848484848485
$ws.Range("E4").Value = 4
$ws.Range("G4").Value = 0.09523809523809523
$ws.Range("H4").Value = 0.8
$ws.Range("J4").Value = 0.9047619047619048
$ws.Range("K4").Value = 31000
$ws.Range("L4").Value = 33
$ws.Range("M4").Value = 5
$ws.Range("N4").Value = 5
$ws.Range("P4").Value = 0.15151515151515152
$ws.Range("Q4").Value = 1
$ws.Range("S4").Value = 0.8484848484848485
$ws.Range("T4").Value = 38600
$ws.Range("U4").Value = 38
$ws.Range("Y4").Value = 0.13157894736842105

# Row 5
$ws.Range("AB5").Value = 0.8
$ws.Range("AC5").Value = 18300
$ws.Range("AD5").Value = 13
$ws.Range("AE5").Value = 5
$ws.Range("AF5").Value = 5
$ws.Range("AH5").Value = 0.38461538461538464
$ws.Range("AI5").Value = 1
$ws.Range("AK5").Value = 0.6153846153846154
$ws.Range("AL5").Value = 16300
$ws.Range("AM5").Value = 15
$ws.Range("AN5").Value = 3
$ws.Range("AO5").Value = 3
$ws.Range("AQ5").Value = 0.2
$ws.Range("AR5").Value = 0.6
$ws.Range("AT5").Value = 0.8
$ws.Range("AU5").Value = 13700
$ws.Range("AV5").Value = 16
$ws.Range("AW5").Value = 4
$ws.Range("AX5").Value = 4
$ws.Range("AZ5").Value = 0.25
$ws.Range("B5").Value = 14300
$ws.Range("BA5").Value = 0.8
$ws.Range("BC5").Value = 0.75
$ws.Range("BD5").Value = 13800
$ws.Range("BE5").Value = 20
$ws.Range("BF5").Value = 4
$ws.Range("BG5").Value = 4
$ws.Range("BI5").Value = 0.2
$ws.Range("BJ5").Value = 0.8
$ws.Range("BL5").Value = 0.8
$ws.Range("BM5").Value = 13400
$ws.Range("BN5").Value = 16
$ws.Range("BO5").Value = 5
$ws.Range("BP5").Value = 5
$ws.Range("BR5").Value = 0.3125
$ws.Range("BS5").Value = 1
$ws.Range("BU5").Value = 0.6875
$ws.Range("BV5").Value = 14900
$ws.Range("BW5").Value = 10
$ws.Range("BX5").Value = 4
$ws.Range("BY5").Value = 4
$ws.Range("C5").Value = 16
$ws.Range("CA5").Value = 0.4
$ws.Range("CB5").Value = 0.8
$ws.Range("CD5").Value = 0.6
$ws.Range("CE5").Value = 15300
$ws.Range("CF5").Value = 17
$ws.Range("CG5").Value = 4
$ws.Range("CH5").Value = 4
$ws.Range("CJ5").Value = 0.23529411764705882
$ws.Range("CK5").Value = 0.8
$ws.Range("CM5").Value = 0.7647058823529411
$ws.Range("CN5").Value = 15180
$ws.Range("CO5").Value = 15.6
$ws.Range("CP5").Value = 4.1
$ws.Range("CQ5").Value = 4.1
$ws.Range("CS5").Value = 0.2740101809954751
$ws.Range("CT5").Value = 0.82
$ws.Range("CV5").Value = 0.7259898190045249
$ws.Range("CW5").Value = 13400
$ws.Range("CX5").Value = 20
$ws.Range("CY5").Value = 5
$ws.Range("CZ5").Value = 5
$ws.Range("D5").Value = 4
$ws.Range("DB5").Value = 0.4
$ws.Range("DC5").Value = 1
$ws.Range("DE5").Value = 0.6153846153846154
$ws.Range("E5").Value = 4
$ws.Range("G5").Value = 0.25
$ws.Range("H5").Value = 0.8
$ws.Range("J5").Value = 0.75
$ws.Range("K5").Value = 15700
$ws.Range("L5").Value = 13
$ws.Range("M5").Value = 4
$ws.Range("N5").Value = 4
$ws.Range("P5").Value = 0.3076923076923077
$ws.Range("Q5").Value = 0.8
$ws.Range("S5").Value = 0.6923076923076923
$ws.Range("T5").Value = 16100
$ws.Range("U5").Value = 20
$ws.Range("V5").Value = 4
$ws.Range("W5").Value = 4
$ws.Range("Y5").Value = 0.2
$ws.Range("Z5").Value = 0.8

# Row 6
$ws.Range("AB6").Value = 0.8979591836734694
$ws.Range("AC6").Value = 52800
$ws.Range("AD6").Value = 55
$ws.Range("AH6").Value = 0.09090909090909091
$ws.Range("AK6").Value = 0.9090909090909091
$ws.Range("AL6").Value = 6100
$ws.Range("AM6").Value = 60
$ws.Range("AQ6").Value = 0.08333333333333333
$ws.Range("AT6").Value = 0.9166666666666666
$ws.Range("AU6").Value = 16700
$ws.Range("AV6").Value = 55
$ws.Range("AZ6").Value = 0.09090909090909091
$ws.Range("B6").Value = 21900
$ws.Range("BC6").Value = 0.9090909090909091
$ws.Range("BD6").Value = 20900
$ws.Range("BE6").Value = 58
$ws.Range("BI6").Value = 0.08620689655172414
$ws.Range("BL6").Value = 0.9137931034482759
$ws.Range("BM6").Value = 35200
$ws.Range("BN6").Value = 59
$ws.Range("BR6").Value = 0.0847457627118644
$ws.Range("BU6").Value = 0.9152542372881356
$ws.Range("BV6").Value = 12400
$ws.Range("BW6").Value = 57
$ws.Range("C6").Value = 58
$ws.Range("CA6").Value = 0.08771929824561403
$ws.Range("CD6").Value = 0.9122807017543859
$ws.Range("CE6").Value = 26900
$ws.Range("CF6").Value = 58
$ws.Range("CG6").Value = 5
$ws.Range("CH6").Value = 5
$ws.Range("CJ6").Value = 0.08620689655172414
$ws.Range("CK6").Value = 1
$ws.Range("CM6").Value = 0.9137931034482759
$ws.Range("CN6").Value = 24800
$ws.Range("CO6").Value = 56.4
$ws.Range("CP6").Value = 5
$ws.Range("CQ6").Value = 5
$ws.Range("CS6").Value = 0.08891871729997876
$ws.Range("CT6").Value = 1
$ws.Range("CV6").Value = 0.9110812827000213
$ws.Range("CW6").Value = 6100
$ws.Range("CX6").Value = 60
$ws.Range("DB6").Value = 0.10204081632653061
$ws.Range("DE6").Value = 0.8979591836734694
$ws.Range("G6").Value = 0.08620689655172414
$ws.Range("J6").Value = 0.9137931034482759
$ws.Range("K6").Value = 43900
$ws.Range("L6").Value = 55
$ws.Range("M6").Value = 5
$ws.Range("N6").Value = 5
$ws.Range("P6").Value = 0.09090909090909091
$ws.Range("Q6").Value = 1
$ws.Range("S6").Value = 0.9090909090909091
$ws.Range("T6").Value = 11200
$ws.Range("U6").Value = 49
$ws.Range("Y6").Value = 0.10204081632653061
